# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" sheets to match the newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5746
$ws1.Range("F4").Value = 98
$ws1.Range("F5").Value = 412
$ws1.Range("F7").Value = 8
$ws1.Range("F8").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5746
$ws4.Range("F4").Value = 98
$ws4.Range("F6").Value = 412
$ws4.Range("F8").Value = 8
$ws4.Range("F9").Value = 24
